$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2025-01-26 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-27 Monday", 1) | Out-Null

# Update each answer cell in the table by position, to avoid ambiguity
# when the same expression text occurs in multiple cells.
$t = $d.Tables(1)

$c = $t.Cell(1, 1)
$c.Range.Find.Execute("33-19=14", $true, $false, $false, $false, $false, $true, 1, $false, "70+25=95", 1) | Out-Null
$c = $t.Cell(1, 2)
$c.Range.Find.Execute("2+13=15", $true, $false, $false, $false, $false, $true, 1, $false, "2+86=88", 1) | Out-Null
$c = $t.Cell(1, 3)
$c.Range.Find.Execute("86-15=71", $true, $false, $false, $false, $false, $true, 1, $false, "52+32=84", 1) | Out-Null
$c = $t.Cell(1, 4)
$c.Range.Find.Execute("33+30=63", $true, $false, $false, $false, $false, $true, 1, $false, "90-53=37", 1) | Out-Null
$c = $t.Cell(1, 5)
$c.Range.Find.Execute("8+72=80", $true, $false, $false, $false, $false, $true, 1, $false, "49+34=83", 1) | Out-Null

$c = $t.Cell(2, 1)
$c.Range.Find.Execute("42-33=9", $true, $false, $false, $false, $false, $true, 1, $false, "21-3=18", 1) | Out-Null
$c = $t.Cell(2, 2)
$c.Range.Find.Execute("15-0=15", $true, $false, $false, $false, $false, $true, 1, $false, "13+31=44", 1) | Out-Null
$c = $t.Cell(2, 3)
$c.Range.Find.Execute("13+57=70", $true, $false, $false, $false, $false, $true, 1, $false, "49+11=60", 1) | Out-Null
$c = $t.Cell(2, 4)
$c.Range.Find.Execute("11+21=32", $true, $false, $false, $false, $false, $true, 1, $false, "8+73=81", 1) | Out-Null
$c = $t.Cell(2, 5)
$c.Range.Find.Execute("85+5=90", $true, $false, $false, $false, $false, $true, 1, $false, "82-56=26", 1) | Out-Null

$c = $t.Cell(3, 1)
$c.Range.Find.Execute("92-54=38", $true, $false, $false, $false, $false, $true, 1, $false, "78+15=93", 1) | Out-Null
$c = $t.Cell(3, 2)
$c.Range.Find.Execute("32-17=15", $true, $false, $false, $false, $false, $true, 1, $false, "45-35=10", 1) | Out-Null
$c = $t.Cell(3, 3)
$c.Range.Find.Execute("19+67=86", $true, $false, $false, $false, $false, $true, 1, $false, "97-36=61", 1) | Out-Null
$c = $t.Cell(3, 4)
$c.Range.Find.Execute("0+42=42", $true, $false, $false, $false, $false, $true, 1, $false, "63-11=52", 1) | Out-Null
$c = $t.Cell(3, 5)
$c.Range.Find.Execute("12+69=81", $true, $false, $false, $false, $false, $true, 1, $false, "57-8=49", 1) | Out-Null

$c = $t.Cell(4, 1)
$c.Range.Find.Execute("31+37=68", $true, $false, $false, $false, $false, $true, 1, $false, "81+1=82", 1) | Out-Null
$c = $t.Cell(4, 2)
$c.Range.Find.Execute("52+47=99", $true, $false, $false, $false, $false, $true, 1, $false, "82-43=39", 1) | Out-Null
$c = $t.Cell(4, 3)
$c.Range.Find.Execute("60-42=18", $true, $false, $false, $false, $false, $true, 1, $false, "51-37=14", 1) | Out-Null
$c = $t.Cell(4, 4)
$c.Range.Find.Execute("76-72=4", $true, $false, $false, $false, $false, $true, 1, $false, "24+37=61", 1) | Out-Null
$c = $t.Cell(4, 5)
$c.Range.Find.Execute("90-57=33", $true, $false, $false, $false, $false, $true, 1, $false, "55+38=93", 1) | Out-Null

$c = $t.Cell(5, 1)
$c.Range.Find.Execute("80-48=32", $true, $false, $false, $false, $false, $true, 1, $false, "50+0=50", 1) | Out-Null
$c = $t.Cell(5, 2)
$c.Range.Find.Execute("61+14=75", $true, $false, $false, $false, $false, $true, 1, $false, "2-1=1", 1) | Out-Null
$c = $t.Cell(5, 3)
$c.Range.Find.Execute("37-5=32", $true, $false, $false, $false, $false, $true, 1, $false, "87+7=94", 1) | Out-Null
$c = $t.Cell(5, 4)
$c.Range.Find.Execute("89-29=60", $true, $false, $false, $false, $false, $true, 1, $false, "34+28=62", 1) | Out-Null
$c = $t.Cell(5, 5)
$c.Range.Find.Execute("8+38=46", $true, $false, $false, $false, $false, $true, 1, $false, "24-18=6", 1) | Out-Null

$c = $t.Cell(6, 1)
$c.Range.Find.Execute("61-56=5", $true, $false, $false, $false, $false, $true, 1, $false, "61-17=44", 1) | Out-Null
$c = $t.Cell(6, 2)
$c.Range.Find.Execute("35+30=65", $true, $false, $false, $false, $false, $true, 1, $false, "75-10=65", 1) | Out-Null
$c = $t.Cell(6, 3)
$c.Range.Find.Execute("24-2=22", $true, $false, $false, $false, $false, $true, 1, $false, "42-3=39", 1) | Out-Null
$c = $t.Cell(6, 4)
$c.Range.Find.Execute("52-25=27", $true, $false, $false, $false, $false, $true, 1, $false, "85+0=85", 1) | Out-Null
$c = $t.Cell(6, 5)
$c.Range.Find.Execute("98-87=11", $true, $false, $false, $false, $false, $true, 1, $false, "40+42=82", 1) | Out-Null

$c = $t.Cell(7, 1)
$c.Range.Find.Execute("82+2=84", $true, $false, $false, $false, $false, $true, 1, $false, "87-1=86", 1) | Out-Null
$c = $t.Cell(7, 2)
$c.Range.Find.Execute("28+4=32", $true, $false, $false, $false, $false, $true, 1, $false, "92+1=93", 1) | Out-Null
$c = $t.Cell(7, 3)
$c.Range.Find.Execute("25+27=52", $true, $false, $false, $false, $false, $true, 1, $false, "41+43=84", 1) | Out-Null
$c = $t.Cell(7, 4)
$c.Range.Find.Execute("75+0=75", $true, $false, $false, $false, $false, $true, 1, $false, "9+24=33", 1) | Out-Null
$c = $t.Cell(7, 5)
$c.Range.Find.Execute("87+0=87", $true, $false, $false, $false, $false, $true, 1, $false, "16+82=98", 1) | Out-Null

$c = $t.Cell(8, 1)
$c.Range.Find.Execute("87-48=39", $true, $false, $false, $false, $false, $true, 1, $false, "15+33=48", 1) | Out-Null
$c = $t.Cell(8, 2)
$c.Range.Find.Execute("54+22=76", $true, $false, $false, $false, $false, $true, 1, $false, "34-12=22", 1) | Out-Null
$c = $t.Cell(8, 3)
$c.Range.Find.Execute("76-74=2", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=71", 1) | Out-Null
$c = $t.Cell(8, 4)
$c.Range.Find.Execute("84+1=85", $true, $false, $false, $false, $false, $true, 1, $false, "2+73=75", 1) | Out-Null
$c = $t.Cell(8, 5)
$c.Range.Find.Execute("75-20=55", $true, $false, $false, $false, $false, $true, 1, $false, "10+54=64", 1) | Out-Null

$c = $t.Cell(9, 1)
$c.Range.Find.Execute("96-86=10", $true, $false, $false, $false, $false, $true, 1, $false, "96-8=88", 1) | Out-Null
$c = $t.Cell(9, 2)
$c.Range.Find.Execute("37+21=58", $true, $false, $false, $false, $false, $true, 1, $false, "26+50=76", 1) | Out-Null
$c = $t.Cell(9, 3)
$c.Range.Find.Execute("90-30=60", $true, $false, $false, $false, $false, $true, 1, $false, "60-7=53", 1) | Out-Null
$c = $t.Cell(9, 4)
$c.Range.Find.Execute("73-31=42", $true, $false, $false, $false, $false, $true, 1, $false, "39-22=17", 1) | Out-Null
$c = $t.Cell(9, 5)
$c.Range.Find.Execute("50+11=61", $true, $false, $false, $false, $false, $true, 1, $false, "50-49=1", 1) | Out-Null

$c = $t.Cell(10, 1)
$c.Range.Find.Execute("34-1=33", $true, $false, $false, $false, $false, $true, 1, $false, "60-18=42", 1) | Out-Null
$c = $t.Cell(10, 2)
$c.Range.Find.Execute("83-32=51", $true, $false, $false, $false, $false, $true, 1, $false, "76-47=29", 1) | Out-Null
$c = $t.Cell(10, 3)
$c.Range.Find.Execute("24+74=98", $true, $false, $false, $false, $false, $true, 1, $false, "14+8=22", 1) | Out-Null
$c = $t.Cell(10, 4)
$c.Range.Find.Execute("20+30=50", $true, $false, $false, $false, $false, $true, 1, $false, "67-20=47", 1) | Out-Null
$c = $t.Cell(10, 5)
$c.Range.Find.Execute("84-67=17", $true, $false, $false, $false, $false, $true, 1, $false, "83-50=33", 1) | Out-Null

$c = $t.Cell(11, 1)
$c.Range.Find.Execute("72-31=41", $true, $false, $false, $false, $false, $true, 1, $false, "35+9=44", 1) | Out-Null
$c = $t.Cell(11, 2)
$c.Range.Find.Execute("40+7=47", $true, $false, $false, $false, $false, $true, 1, $false, "75-58=17", 1) | Out-Null
$c = $t.Cell(11, 3)
$c.Range.Find.Execute("66-13=53", $true, $false, $false, $false, $false, $true, 1, $false, "98-96=2", 1) | Out-Null
$c = $t.Cell(11, 4)
$c.Range.Find.Execute("89-28=61", $true, $false, $false, $false, $false, $true, 1, $false, "56-42=14", 1) | Out-Null
$c = $t.Cell(11, 5)
$c.Range.Find.Execute("85-74=11", $true, $false, $false, $false, $false, $true, 1, $false, "14+75=89", 1) | Out-Null

$c = $t.Cell(12, 1)
$c.Range.Find.Execute("82-27=55", $true, $false, $false, $false, $false, $true, 1, $false, "28+55=83", 1) | Out-Null
$c = $t.Cell(12, 2)
$c.Range.Find.Execute("38+17=55", $true, $false, $false, $false, $false, $true, 1, $false, "11+34=45", 1) | Out-Null
$c = $t.Cell(12, 3)
$c.Range.Find.Execute("56-51=5", $true, $false, $false, $false, $false, $true, 1, $false, "54+33=87", 1) | Out-Null
$c = $t.Cell(12, 4)
$c.Range.Find.Execute("48-32=16", $true, $false, $false, $false, $false, $true, 1, $false, "69-41=28", 1) | Out-Null
$c = $t.Cell(12, 5)
$c.Range.Find.Execute("47+35=82", $true, $false, $false, $false, $false, $true, 1, $false, "9+65=74", 1) | Out-Null

$c = $t.Cell(13, 1)
$c.Range.Find.Execute("82-2=80", $true, $false, $false, $false, $false, $true, 1, $false, "51-0=51", 1) | Out-Null
$c = $t.Cell(13, 2)
$c.Range.Find.Execute("27+4=31", $true, $false, $false, $false, $false, $true, 1, $false, "49-40=9", 1) | Out-Null
$c = $t.Cell(13, 3)
$c.Range.Find.Execute("58-5=53", $true, $false, $false, $false, $false, $true, 1, $false, "42+1=43", 1) | Out-Null
$c = $t.Cell(13, 4)
$c.Range.Find.Execute("30+28=58", $true, $false, $false, $false, $false, $true, 1, $false, "15-2=13", 1) | Out-Null
$c = $t.Cell(13, 5)
$c.Range.Find.Execute("79-55=24", $true, $false, $false, $false, $false, $true, 1, $false, "5+22=27", 1) | Out-Null

$c = $t.Cell(14, 1)
$c.Range.Find.Execute("13+2=15", $true, $false, $false, $false, $false, $true, 1, $false, "99-24=75", 1) | Out-Null
$c = $t.Cell(14, 2)
$c.Range.Find.Execute("74-55=19", $true, $false, $false, $false, $false, $true, 1, $false, "64-21=43", 1) | Out-Null
$c = $t.Cell(14, 3)
$c.Range.Find.Execute("52-15=37", $true, $false, $false, $false, $false, $true, 1, $false, "51-43=8", 1) | Out-Null
$c = $t.Cell(14, 4)
$c.Range.Find.Execute("93-29=64", $true, $false, $false, $false, $false, $true, 1, $false, "84-32=52", 1) | Out-Null
$c = $t.Cell(14, 5)
$c.Range.Find.Execute("5+94=99", $true, $false, $false, $false, $false, $true, 1, $false, "36+55=91", 1) | Out-Null

$c = $t.Cell(15, 1)
$c.Range.Find.Execute("87-75=12", $true, $false, $false, $false, $false, $true, 1, $false, "80-28=52", 1) | Out-Null
$c = $t.Cell(15, 2)
$c.Range.Find.Execute("21+77=98", $true, $false, $false, $false, $false, $true, 1, $false, "99-12=87", 1) | Out-Null
$c = $t.Cell(15, 3)
$c.Range.Find.Execute("79-52=27", $true, $false, $false, $false, $false, $true, 1, $false, "86-46=40", 1) | Out-Null
$c = $t.Cell(15, 4)
$c.Range.Find.Execute("29-10=19", $true, $false, $false, $false, $false, $true, 1, $false, "24-23=1", 1) | Out-Null
$c = $t.Cell(15, 5)
$c.Range.Find.Execute("26+23=49", $true, $false, $false, $false, $false, $true, 1, $false, "37+6=43", 1) | Out-Null

$c = $t.Cell(16, 1)
$c.Range.Find.Execute("43-39=4", $true, $false, $false, $false, $false, $true, 1, $false, "3+27=30", 1) | Out-Null
$c = $t.Cell(16, 2)
$c.Range.Find.Execute("84-53=31", $true, $false, $false, $false, $false, $true, 1, $false, "47+22=69", 1) | Out-Null
$c = $t.Cell(16, 3)
$c.Range.Find.Execute("49-45=4", $true, $false, $false, $false, $false, $true, 1, $false, "44+20=64", 1) | Out-Null
$c = $t.Cell(16, 4)
$c.Range.Find.Execute("32-13=19", $true, $false, $false, $false, $false, $true, 1, $false, "14+39=53", 1) | Out-Null
$c = $t.Cell(16, 5)
$c.Range.Find.Execute("1+91=92", $true, $false, $false, $false, $false, $true, 1, $false, "41-23=18", 1) | Out-Null

$c = $t.Cell(17, 1)
$c.Range.Find.Execute("4+6=10", $true, $false, $false, $false, $false, $true, 1, $false, "12+17=29", 1) | Out-Null
$c = $t.Cell(17, 2)
$c.Range.Find.Execute("6+45=51", $true, $false, $false, $false, $false, $true, 1, $false, "69-14=55", 1) | Out-Null
$c = $t.Cell(17, 3)
$c.Range.Find.Execute("60-10=50", $true, $false, $false, $false, $false, $true, 1, $false, "8+77=85", 1) | Out-Null
$c = $t.Cell(17, 4)
$c.Range.Find.Execute("93-58=35", $true, $false, $false, $false, $false, $true, 1, $false, "83-22=61", 1) | Out-Null
$c = $t.Cell(17, 5)
$c.Range.Find.Execute("33-15=18", $true, $false, $false, $false, $false, $true, 1, $false, "61-55=6", 1) | Out-Null

$c = $t.Cell(18, 1)
$c.Range.Find.Execute("80-3=77", $true, $false, $false, $false, $false, $true, 1, $false, "83-49=34", 1) | Out-Null
$c = $t.Cell(18, 2)
$c.Range.Find.Execute("60-42=18", $true, $false, $false, $false, $false, $true, 1, $false, "14+63=77", 1) | Out-Null
$c = $t.Cell(18, 3)
$c.Range.Find.Execute("23-5=18", $true, $false, $false, $false, $false, $true, 1, $false, "75-2=73", 1) | Out-Null
$c = $t.Cell(18, 4)
$c.Range.Find.Execute("36-14=22", $true, $false, $false, $false, $false, $true, 1, $false, "57-15=42", 1) | Out-Null
$c = $t.Cell(18, 5)
$c.Range.Find.Execute("48+47=95", $true, $false, $false, $false, $false, $true, 1, $false, "11+63=74", 1) | Out-Null

$c = $t.Cell(19, 1)
$c.Range.Find.Execute("67-34=33", $true, $false, $false, $false, $false, $true, 1, $false, "29+57=86", 1) | Out-Null
$c = $t.Cell(19, 2)
$c.Range.Find.Execute("26+62=88", $true, $false, $false, $false, $false, $true, 1, $false, "0+20=20", 1) | Out-Null
$c = $t.Cell(19, 3)
$c.Range.Find.Execute("84-82=2", $true, $false, $false, $false, $false, $true, 1, $false, "15+23=38", 1) | Out-Null
$c = $t.Cell(19, 4)
$c.Range.Find.Execute("30-11=19", $true, $false, $false, $false, $false, $true, 1, $false, "84-11=73", 1) | Out-Null
$c = $t.Cell(19, 5)
$c.Range.Find.Execute("25+21=46", $true, $false, $false, $false, $false, $true, 1, $false, "88-36=52", 1) | Out-Null

$c = $t.Cell(20, 1)
$c.Range.Find.Execute("69+1=70", $true, $false, $false, $false, $false, $true, 1, $false, "41+36=77", 1) | Out-Null
$c = $t.Cell(20, 2)
$c.Range.Find.Execute("98-56=42", $true, $false, $false, $false, $false, $true, 1, $false, "84-54=30", 1) | Out-Null
$c = $t.Cell(20, 3)
$c.Range.Find.Execute("56+8=64", $true, $false, $false, $false, $false, $true, 1, $false, "20+21=41", 1) | Out-Null
$c = $t.Cell(20, 4)
$c.Range.Find.Execute("36+20=56", $true, $false, $false, $false, $false, $true, 1, $false, "92-3=89", 1) | Out-Null
$c = $t.Cell(20, 5)
$c.Range.Find.Execute("34+5=39", $true, $false, $false, $false, $false, $true, 1, $false, "77-13=64", 1) | Out-Null
